# Updates Bdnf-Ntrk2 LR-pair sheet with refreshed TPM-based values and adds the
# new "Resolving-Mac" target cluster rows (8 and 9).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = @("FAPs", "Bdnf", "Ntrk2", "ECs", 1, 0.3333333333333333, 0.3379886666666667, 1.013966, 0.4664623486936776, 0.4664623486936776, 3, 1, 0.293221, 0.879663, 0.02611983441994871, 0.02611983441994871, 0.09910537482866667, 0.8919483734579999, 0.01218391931101924, 0.01218391931101924)
$row3 = @("FAPs", "Bdnf", "Ntrk2", "FAPs", 1, 0.3333333333333333, 0.3379886666666667, 1.013966, 0.4664623486936776, 0.4664623486936776, 3, 1, 10.74310933333333, 32.229328, 0.9569854715114954, 0.9569854715114955, 3.631049199427555, 32.67944279484799, 0.4463976907069786, 0.4463976907069787)
$row4 = @("FAPs", "Bdnf", "Ntrk2", "MuSCs", 1, 0.3333333333333333, 0.3379886666666667, 1.013966, 0.4664623486936776, 0.4664623486936776, 3, 1, 0.183427, 0.550281, 0.01633949433413, 0.01633949433413, 0.06199624716066666, 0.557966224446, 0.007621758903565316, 0.007621758903565317)
$row5 = @("FAPs", "Bdnf", "Ntrk2", "Resolving-Mac", 1, 0.3333333333333333, 0.3379886666666667, 1.013966, 0.4664623486936776, 0.4664623486936776, 1, 0.3333333333333333, 0.006232666666666667, 0.018698, 0.0005551997344257983, 0.0005551997344257983, 0.002106570696444445, 0.018959136268, 0.0002589797721143639, 0.0002589797721143639)
$row6 = @("MuSCs", "Bdnf", "Ntrk2", "ECs", 2, 0.6666666666666666, 0.38659, 1.15977, 0.5335376513063224, 0.5335376513063224, 3, 1, 0.293221, 0.879663, 0.02611983441994871, 0.02611983441994871, 0.11335630639, 1.02020675751, 0.01393591510892947, 0.01393591510892947)
$row7 = @("MuSCs", "Bdnf", "Ntrk2", "FAPs", 2, 0.6666666666666666, 0.38659, 1.15977, 0.5335376513063224, 0.5335376513063224, 3, 1, 10.74310933333333, 32.229328, 0.9569854715114954, 0.9569854715114955, 4.153178637173332, 37.37860773455999, 0.5105877808045167, 0.5105877808045167)
$row8 = @("MuSCs", "Bdnf", "Ntrk2", "MuSCs", 2, 0.6666666666666666, 0.38659, 1.15977, 0.5335376513063224, 0.5335376513063224, 3, 1, 0.183427, 0.550281, 0.01633949433413, 0.01633949433413, 0.07091104393, 0.63819939537, 0.00871773543056468, 0.008717735430564682)
$row9 = @("MuSCs", "Bdnf", "Ntrk2", "Resolving-Mac", 2, 0.6666666666666666, 0.38659, 1.15977, 0.5335376513063224, 0.5335376513063224, 1, 0.3333333333333333, 0.006232666666666667, 0.018698, 0.0005551997344257983, 0.0005551997344257983, 0.002409486606666667, 0.02168537946, 0.0002962199623114343, 0.0002962199623114343)

$rows = @(2, 3, 4, 5, 6, 7, 8, 9)
$data = @($row2, $row3, $row4, $row5, $row6, $row7, $row8, $row9)

for ($i = 0; $i -lt $rows.Count; $i++) {
  $r = $rows[$i]
  $vals = $data[$i]
  for ($c = 1; $c -le $vals.Count; $c++) {
    $ws.Cells.Item($r, $c).Value = $vals[$c - 1]
  }
}
